# Refresh the "想去人数" (want-to-go count) figures in column F across all
# four sheets, matching the regenerated scrape output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 359
$ws.Range("F5").Value = 8391
$ws.Range("F8").Value = 2280
$ws.Range("F13").Value = 7617
$ws.Range("F14").Value = 7829
$ws.Range("F15").Value = 58320
$ws.Range("F16").Value = 4981
$ws.Range("F23").Value = 5356
$ws.Range("F29").Value = 1466
$ws.Range("F33").Value = 256
$ws.Range("F36").Value = 746
$ws.Range("F38").Value = 797
$ws.Range("F39").Value = 1217
$ws.Range("F40").Value = 427
$ws.Range("F41").Value = 17
$ws.Range("F43").Value = 231

# Sheet 2: 演出 (Performances)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 180
$ws.Range("F5").Value = 7755

# Sheet 3: 本地生活 (Local life)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F7").Value = 691
$ws.Range("F16").Value = 2571
$ws.Range("F17").Value = 273

# Sheet 4: 全部类型 (All types — combined listing)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 359
$ws.Range("F5").Value = 691
$ws.Range("F8").Value = 2571
$ws.Range("F9").Value = 273
$ws.Range("F11").Value = 7617
$ws.Range("F12").Value = 7829
$ws.Range("F13").Value = 4981
$ws.Range("F17").Value = 180
$ws.Range("F18").Value = 5356
$ws.Range("F22").Value = 1466
$ws.Range("F32").Value = 746
$ws.Range("F34").Value = 797
$ws.Range("F38").Value = 427
$ws.Range("F42").Value = 231
